$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Cells.Item(12, 8).Value = 0  # H12: 195 -> 0
$ws.Cells.Item(12, 9).Value = 0  # I12: 195 -> 0
$ws.Cells.Item(12, 11).Value = 0  # K12: 195 -> 0
$ws.Cells.Item(12, 13).ClearContents()  # M12: -25 -> (removed)
# Row 17
$ws.Cells.Item(17, 8).Value = 1999  # H17: 0 -> 1999
$ws.Cells.Item(17, 10).Value = 1999  # J17: 0 -> 1999
$ws.Cells.Item(17, 12).Value = 5997  # L17: 0 -> 5997
$ws.Cells.Item(17, 14).Value = -6333  # N17: None -> -6333
# Row 26
$ws.Cells.Item(26, 8).Value = 501  # H26: 0 -> 501
$ws.Cells.Item(26, 9).Value = 501  # I26: 0 -> 501
$ws.Cells.Item(26, 11).Value = 501  # K26: 0 -> 501
$ws.Cells.Item(26, 13).Value = -157  # M26: None -> -157
# Row 41
$ws.Cells.Item(41, 8).Value = 2211.7273  # H41: 2285.7273 -> 2211.7273
$ws.Cells.Item(41, 9).Value = 121.666664  # I41: 252.6 -> 121.666664
$ws.Cells.Item(41, 10).Value = 2995.5  # J41: 3980 -> 2995.5
$ws.Cells.Item(41, 11).Value = 121.666664  # K41: 252.6 -> 121.666664
$ws.Cells.Item(41, 12).Value = 2995.5  # L41: 3980 -> 2995.5
$ws.Cells.Item(41, 13).Value = 318.333336  # M41: 187.4 -> 318.333336
$ws.Cells.Item(41, 14).Value = -3875.5  # N41: -4860 -> -3875.5
# Row 112
$ws.Cells.Item(112, 8).Value = 6999.5  # H112: 5333 -> 6999.5
$ws.Cells.Item(112, 9).Value = 0  # I112: 2000 -> 0
$ws.Cells.Item(112, 11).Value = 0  # K112: 6000 -> 0
$ws.Cells.Item(112, 13).ClearContents()  # M112: -4892 -> (removed)
# Row 132
$ws.Cells.Item(132, 8).Value = 3049.6191  # H132: 3184.8 -> 3049.6191
$ws.Cells.Item(132, 9).Value = 1106  # I132: 1156.4286 -> 1106
$ws.Cells.Item(132, 10).Value = 6936.857  # J132: 7917.6665 -> 6936.857
$ws.Cells.Item(132, 11).Value = 3318  # K132: 3469.2858 -> 3318
$ws.Cells.Item(132, 12).Value = 20810.571  # L132: 23752.9995 -> 20810.571
$ws.Cells.Item(132, 13).Value = -788  # M132: -939.2857999999997 -> -788
$ws.Cells.Item(132, 14).Value = -25870.571  # N132: -28812.9995 -> -25870.571
# Row 141
$ws.Cells.Item(141, 8).Value = 40099  # H141: 49998.75 -> 40099
$ws.Cells.Item(141, 9).Value = 40099  # I141: 49998.75 -> 40099
$ws.Cells.Item(141, 11).Value = 120297  # K141: 149996.25 -> 120297
$ws.Cells.Item(141, 13).Value = -115117  # M141: -144816.25 -> -115117

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 2293.5  # H45: 2298.5 -> 2293.5
$ws.Cells.Item(45, 10).Value = 1649.5  # J45: 1699.5 -> 1649.5
$ws.Cells.Item(45, 12).Value = 1649.5  # L45: 1699.5 -> 1649.5
$ws.Cells.Item(45, 14).Value = -2403.5  # N45: -2453.5 -> -2403.5
# Row 61
$ws.Cells.Item(61, 8).Value = 1000  # H61: 0 -> 1000
$ws.Cells.Item(61, 9).Value = 1000  # I61: 0 -> 1000
$ws.Cells.Item(61, 11).Value = 1000  # K61: 0 -> 1000
$ws.Cells.Item(61, 13).Value = -788  # M61: None -> -788
# Row 132
$ws.Cells.Item(132, 8).Value = 2200  # H132: 2240 -> 2200
$ws.Cells.Item(132, 9).Value = 2200  # I132: 2240 -> 2200
$ws.Cells.Item(132, 11).Value = 6600  # K132: 6720 -> 6600
$ws.Cells.Item(132, 13).Value = -4070  # M132: -4190 -> -4070
# Row 136
$ws.Cells.Item(136, 8).Value = 1000  # H136: 0 -> 1000
$ws.Cells.Item(136, 9).Value = 1000  # I136: 0 -> 1000
$ws.Cells.Item(136, 11).Value = 3000  # K136: 0 -> 3000
$ws.Cells.Item(136, 13).Value = -450  # M136: None -> -450

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1502.25  # H20: 2061.6 -> 1502.25
$ws.Cells.Item(20, 9).Value = 1375  # I20: 2250 -> 1375
$ws.Cells.Item(20, 10).Value = 1629.5  # J20: 1936 -> 1629.5
$ws.Cells.Item(20, 11).Value = 1375  # K20: 2250 -> 1375
$ws.Cells.Item(20, 12).Value = 1629.5  # L20: 1936 -> 1629.5
$ws.Cells.Item(20, 13).Value = -1128  # M20: -2003 -> -1128
$ws.Cells.Item(20, 14).Value = -2123.5  # N20: -2430 -> -2123.5
# Row 22
$ws.Cells.Item(22, 8).Value = 3370853.2  # H22: 5056229 -> 3370853.2
$ws.Cells.Item(22, 9).Value = 3370853.2  # I22: 5056229 -> 3370853.2
$ws.Cells.Item(22, 11).Value = 3370853.2  # K22: 5056229 -> 3370853.2
$ws.Cells.Item(22, 13).Value = -3370680.2  # M22: -5056056 -> -3370680.2
# Row 107
$ws.Cells.Item(107, 8).Value = 18625.715  # H107: 10416.77 -> 18625.715
$ws.Cells.Item(107, 9).Value = 8796.4  # I107: 4456.364 -> 8796.4
$ws.Cells.Item(107, 11).Value = 8796.4  # K107: 4456.364 -> 8796.4
$ws.Cells.Item(107, 13).Value = -6876.4  # M107: -2536.364 -> -6876.4
# Row 134
$ws.Cells.Item(134, 8).Value = 1285.6  # H134: 1476.2222 -> 1285.6
$ws.Cells.Item(134, 9).Value = 1132  # I134: 1414.3334 -> 1132
$ws.Cells.Item(134, 10).Value = 1900  # J134: 1600 -> 1900
$ws.Cells.Item(134, 11).Value = 3396  # K134: 4243.0002 -> 3396
$ws.Cells.Item(134, 12).Value = 5700  # L134: 4800 -> 5700
$ws.Cells.Item(134, 13).Value = -861  # M134: -1708.0002 -> -861
$ws.Cells.Item(134, 14).Value = -10770  # N134: -9870 -> -10770

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Cells.Item(23, 8).Value = 170000  # H23: 0 -> 170000
$ws.Cells.Item(23, 9).Value = 170000  # I23: 0 -> 170000
$ws.Cells.Item(23, 11).Value = 170000  # K23: 0 -> 170000
$ws.Cells.Item(23, 13).Value = -169760  # M23: None -> -169760
# Row 27
$ws.Cells.Item(27, 8).Value = 170000  # H27: 0 -> 170000
$ws.Cells.Item(27, 9).Value = 170000  # I27: 0 -> 170000
$ws.Cells.Item(27, 11).Value = 170000  # K27: 0 -> 170000
$ws.Cells.Item(27, 13).Value = -169808  # M27: None -> -169808
# Row 58
$ws.Cells.Item(58, 8).Value = 2510.6667  # H58: 2481.3333 -> 2510.6667
$ws.Cells.Item(58, 9).Value = 2412.8  # I58: 2481.3333 -> 2412.8
$ws.Cells.Item(58, 10).Value = 3000  # J58: 0 -> 3000
$ws.Cells.Item(58, 11).Value = 2412.8  # K58: 2481.3333 -> 2412.8
$ws.Cells.Item(58, 12).Value = 3000  # L58: 0 -> 3000
$ws.Cells.Item(58, 13).Value = -2209.8  # M58: -2278.3333 -> -2209.8
$ws.Cells.Item(58, 14).Value = -3406  # N58: None -> -3406
# Row 69
$ws.Cells.Item(69, 8).Value = 10000  # H69: 8000 -> 10000
$ws.Cells.Item(69, 9).Value = 10000  # I69: 8000 -> 10000
$ws.Cells.Item(69, 11).Value = 10000  # K69: 8000 -> 10000
$ws.Cells.Item(69, 13).Value = -9251  # M69: -7251 -> -9251
# Row 72
$ws.Cells.Item(72, 8).Value = 10000  # H72: 8000 -> 10000
$ws.Cells.Item(72, 9).Value = 10000  # I72: 8000 -> 10000
$ws.Cells.Item(72, 11).Value = 30000  # K72: 24000 -> 30000
$ws.Cells.Item(72, 13).Value = -26256  # M72: -20256 -> -26256
# Row 99
$ws.Cells.Item(99, 8).Value = 2056  # H99: 2075.6 -> 2056
$ws.Cells.Item(99, 9).Value = 2093.3333  # I99: 2126 -> 2093.3333
$ws.Cells.Item(99, 11).Value = 2093.3333  # K99: 2126 -> 2093.3333
$ws.Cells.Item(99, 13).Value = -595.3332999999998  # M99: -628 -> -595.3332999999998
# Row 126
$ws.Cells.Item(126, 8).Value = 2056  # H126: 2075.6 -> 2056
$ws.Cells.Item(126, 9).Value = 2093.3333  # I126: 2126 -> 2093.3333
$ws.Cells.Item(126, 11).Value = 6279.999899999999  # K126: 6378 -> 6279.999899999999
$ws.Cells.Item(126, 13).Value = -3809.999899999999  # M126: -3908 -> -3809.999899999999
# Row 132
$ws.Cells.Item(132, 8).Value = 3800.0476  # H132: 3999.85 -> 3800.0476
$ws.Cells.Item(132, 9).Value = 3438.625  # I132: 3765.2856 -> 3438.625
$ws.Cells.Item(132, 10).Value = 4956.6  # J132: 4547.1665 -> 4956.6
$ws.Cells.Item(132, 11).Value = 10315.875  # K132: 11295.8568 -> 10315.875
$ws.Cells.Item(132, 12).Value = 14869.8  # L132: 13641.4995 -> 14869.8
$ws.Cells.Item(132, 13).Value = -7785.875  # M132: -8765.856800000001 -> -7785.875
$ws.Cells.Item(132, 14).Value = -19929.8  # N132: -18701.4995 -> -19929.8
# Row 136
$ws.Cells.Item(136, 8).Value = 2510.6667  # H136: 2481.3333 -> 2510.6667
$ws.Cells.Item(136, 9).Value = 2412.8  # I136: 2481.3333 -> 2412.8
$ws.Cells.Item(136, 10).Value = 3000  # J136: 0 -> 3000
$ws.Cells.Item(136, 11).Value = 7238.400000000001  # K136: 7443.999899999999 -> 7238.400000000001
$ws.Cells.Item(136, 12).Value = 9000  # L136: 0 -> 9000
$ws.Cells.Item(136, 13).Value = -4688.400000000001  # M136: -4893.999899999999 -> -4688.400000000001
$ws.Cells.Item(136, 14).Value = -14100  # N136: None -> -14100

$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Cells.Item(8, 8).Value = 1291  # H8: 1178.4 -> 1291
$ws.Cells.Item(8, 9).Value = 1291  # I8: 1178.4 -> 1291
$ws.Cells.Item(8, 11).Value = 3873  # K8: 3535.2 -> 3873
$ws.Cells.Item(8, 13).Value = -3734  # M8: -3396.2 -> -3734
# Row 113
$ws.Cells.Item(113, 8).Value = 884  # H113: 1098.25 -> 884
$ws.Cells.Item(113, 9).Value = 441  # I113: 0 -> 441
$ws.Cells.Item(113, 10).Value = 1031.6666  # J113: 1098.25 -> 1031.6666
$ws.Cells.Item(113, 11).Value = 1323  # K113: 0 -> 1323
$ws.Cells.Item(113, 12).Value = 3094.9998  # L113: 3294.75 -> 3094.9998
$ws.Cells.Item(113, 13).Value = 847  # M113: None -> 847
$ws.Cells.Item(113, 14).Value = -7434.9998  # N113: -7634.75 -> -7434.9998
# Row 134
$ws.Cells.Item(134, 8).Value = 1497.25  # H134: 2000 -> 1497.25
$ws.Cells.Item(134, 9).Value = 1497.25  # I134: 2000 -> 1497.25
$ws.Cells.Item(134, 11).Value = 4491.75  # K134: 6000 -> 4491.75
$ws.Cells.Item(134, 13).Value = 578.25  # M134: -930 -> 578.25

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Cells.Item(13, 8).Value = 1351  # H13: 1000 -> 1351
$ws.Cells.Item(13, 9).Value = 0  # I13: 1000 -> 0
$ws.Cells.Item(13, 10).Value = 1351  # J13: 0 -> 1351
$ws.Cells.Item(13, 11).Value = 0  # K13: 1000 -> 0
$ws.Cells.Item(13, 12).Value = 1351  # L13: 0 -> 1351
$ws.Cells.Item(13, 13).ClearContents()  # M13: -861 -> (removed)
$ws.Cells.Item(13, 14).Value = -1629  # N13: None -> -1629
# Row 126
$ws.Cells.Item(126, 8).Value = 4965  # H126: 4975 -> 4965
$ws.Cells.Item(126, 9).Value = 4945  # I126: 0 -> 4945
$ws.Cells.Item(126, 11).Value = 14835  # K126: 0 -> 14835
$ws.Cells.Item(126, 13).Value = -12365  # M126: None -> -12365
# Row 128
$ws.Cells.Item(128, 8).Value = 45000  # H128: 0 -> 45000
$ws.Cells.Item(128, 10).Value = 45000  # J128: 0 -> 45000
$ws.Cells.Item(128, 12).Value = 45000  # L128: 0 -> 45000
$ws.Cells.Item(128, 14).Value = -54960  # N128: None -> -54960
# Row 132
$ws.Cells.Item(132, 8).Value = 2832.5  # H132: 3203.5715 -> 2832.5
$ws.Cells.Item(132, 9).Value = 2980.5557  # I132: 3487.5 -> 2980.5557
$ws.Cells.Item(132, 11).Value = 8941.667099999999  # K132: 10462.5 -> 8941.667099999999
$ws.Cells.Item(132, 13).Value = -6411.667099999999  # M132: -7932.5 -> -6411.667099999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 7166.3335  # H7: 7750 -> 7166.3335
$ws.Cells.Item(7, 9).Value = 6999.5  # I7: 8000 -> 6999.5
$ws.Cells.Item(7, 11).Value = 6999.5  # K7: 8000 -> 6999.5
$ws.Cells.Item(7, 13).Value = -6887.5  # M7: -7888 -> -6887.5
# Row 40
$ws.Cells.Item(40, 8).Value = 3349  # H40: 3389 -> 3349
$ws.Cells.Item(40, 9).Value = 3349  # I40: 3389 -> 3349
$ws.Cells.Item(40, 11).Value = 3349  # K40: 3389 -> 3349
$ws.Cells.Item(40, 13).Value = -3213  # M40: -3253 -> -3213
# Row 126
$ws.Cells.Item(126, 8).Value = 7166.3335  # H126: 7750 -> 7166.3335
$ws.Cells.Item(126, 9).Value = 6999.5  # I126: 8000 -> 6999.5
$ws.Cells.Item(126, 11).Value = 20998.5  # K126: 24000 -> 20998.5
$ws.Cells.Item(126, 13).Value = -18528.5  # M126: -21530 -> -18528.5
# Row 136
$ws.Cells.Item(136, 8).Value = 4943.1763  # H136: 4580.684 -> 4943.1763
$ws.Cells.Item(136, 9).Value = 4814.625  # I136: 4446.278 -> 4814.625
$ws.Cells.Item(136, 11).Value = 14443.875  # K136: 13338.834 -> 14443.875
$ws.Cells.Item(136, 13).Value = -11893.875  # M136: -10788.834 -> -11893.875

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Cells.Item(34, 8).Value = 0  # H34: 17998 -> 0
$ws.Cells.Item(34, 9).Value = 0  # I34: 17998 -> 0
$ws.Cells.Item(34, 11).Value = 0  # K34: 17998 -> 0
$ws.Cells.Item(34, 13).ClearContents()  # M34: -17795 -> (removed)
# Row 100
$ws.Cells.Item(100, 8).Value = 597  # H100: 6154 -> 597
$ws.Cells.Item(100, 9).Value = 484.625  # I100: 7082.1816 -> 484.625
$ws.Cells.Item(100, 10).Value = 1046.5  # J100: 1049 -> 1046.5
$ws.Cells.Item(100, 11).Value = 969.25  # K100: 14164.3632 -> 969.25
$ws.Cells.Item(100, 12).Value = 2093  # L100: 2098 -> 2093
$ws.Cells.Item(100, 13).Value = -428.25  # M100: -13623.3632 -> -428.25
$ws.Cells.Item(100, 14).Value = -3175  # N100: -3180 -> -3175
# Row 126
$ws.Cells.Item(126, 8).Value = 2733  # H126: 2956.6 -> 2733
$ws.Cells.Item(126, 9).Value = 2406.3  # I126: 2618.4443 -> 2406.3
$ws.Cells.Item(126, 11).Value = 7218.900000000001  # K126: 7855.3329 -> 7218.900000000001
$ws.Cells.Item(126, 13).Value = -4748.900000000001  # M126: -5385.3329 -> -4748.900000000001
